# Apply the "additional scraping" update:
#  - insert a new "Player Info" sheet at the front with ID/NAME/BATTING_HAND/BOWL_STYLE
#  - rename the MATCH_CARD_LINK column to MATCH_CODE on both "ODI Batting" and
#    "ODI Bowling" sheets, replacing the full scorecard URL with just the numeric
#    match code that used to be the MatchCode= query parameter.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Insert the new "Player Info" worksheet before "ODI Batting" (i.e. first).
#    NOTE: re-fetch sheet references by name AFTER this insertion, since sheet
#    handles obtained beforehand track position rather than identity.
# ---------------------------------------------------------------------------
$battingSheetBefore = $wb.Worksheets.Item("ODI Batting")
$infoSheet = $wb.Worksheets.Add($battingSheetBefore)
$infoSheet.Name = "Player Info"

$infoSheet.Range("A1").Value = "ID"
$infoSheet.Range("B1").Value = "NAME"
$infoSheet.Range("C1").Value = "BATTING_HAND"
$infoSheet.Range("D1").Value = "BOWL_STYLE"

# Match the bold / bordered / centred-top header styling used on the other sheets.
$infoHeader = $infoSheet.Range("A1:D1")
$infoHeader.Font.Bold = $true
$infoHeader.HorizontalAlignment = -4108
$infoHeader.VerticalAlignment = -4160
$infoHeader.Borders.LineStyle = 1

# Player id needs to stay textual (it's stored as text everywhere else in the workbook).
$infoSheet.Range("A2").NumberFormat = "@"
$infoSheet.Range("A2").Value = "7155"
$infoSheet.Range("B2").Value = "Gerald Coetzee"
$infoSheet.Range("C2").Value = "Right Handed"
$infoSheet.Range("D2").Value = "Right Arm Fast"

# ---------------------------------------------------------------------------
# 2. "ODI Batting": MATCH_CARD_LINK (column D) -> MATCH_CODE, values become the
#    bare numeric match code instead of the full scorecard URL.
# ---------------------------------------------------------------------------
$battingSheet = $wb.Worksheets.Item("ODI Batting")
$battingSheet.Range("D1").Value = "MATCH_CODE"

$battingSheet.Range("D2:D3").NumberFormat = "@"
$battingSheet.Range("D2").Value = "4727"
$battingSheet.Range("D3").Value = "4731"

# ---------------------------------------------------------------------------
# 3. "ODI Bowling": MATCH_CARD_LINK (column B) -> MATCH_CODE, same treatment.
# ---------------------------------------------------------------------------
$bowlingSheet = $wb.Worksheets.Item("ODI Bowling")
$bowlingSheet.Range("B1").Value = "MATCH_CODE"

$bowlingSheet.Range("B2:B3").NumberFormat = "@"
$bowlingSheet.Range("B2").Value = "4727"
$bowlingSheet.Range("B3").Value = "4731"

Write-Output "Sheets now:"
foreach ($s in $wb.Worksheets) {
    Write-Output "  [$($s.Index)] $($s.Name)"
}
